# Fruta / hortaliza, semanal
# Inserts 3 new price rows (week of 2021-10-07, serial 44476) right after the
# existing "44468" block (old row 351) in the Femacal de La Calera - Kiwi
# sheet, pushing every subsequent row down by 3 (old A1:T439 -> A1:T442).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at row 351; Excel shifts rows 351:439 down to 354:442
# and copies formatting (incl. the date-style on column D) from the row above.
$ws.Rows.Item(351).Resize(3).Insert()

# Fixed/common values shared by the three new rows.
$mercadoId = 3
$mercado   = "Femacal de La Calera"
$region    = "Coquimbo"
$codreg    = 5
$tipo      = "Fruta"
$productoId   = 100101
$producto     = "Berries"
$categoriaId  = 100101007
$categoria    = "Kiwi"
$variedad     = "Hayward"
$unidad       = "`$/bandeja 10 kilos"
$origen       = "Región de O'Higgins"
$kgPorUnidad  = 10
$fecha        = 44476

# Per-quality values: Calidad, Volumen, Precio (min=max=prom), Precio/Kg
$newRows = @(
    @{ Row = 351; Calidad = "Especial"; Volumen = 60; Precio = 12000; PrecioKg = 1200 },
    @{ Row = 352; Calidad = "Primera";  Volumen = 68; Precio = 11000; PrecioKg = 1100 },
    @{ Row = 353; Calidad = "Segunda";  Volumen = 65; Precio = 10000; PrecioKg = 1000 }
)

foreach ($nr in $newRows) {
    $r = $nr.Row
    $ws.Cells.Item($r, 1).Value  = $mercadoId
    $ws.Cells.Item($r, 2).Value  = $mercado
    $ws.Cells.Item($r, 3).Value  = $region
    $ws.Cells.Item($r, 4).Value  = $fecha
    $ws.Cells.Item($r, 5).Value  = $codreg
    $ws.Cells.Item($r, 6).Value  = $tipo
    $ws.Cells.Item($r, 7).Value  = $productoId
    $ws.Cells.Item($r, 8).Value  = $producto
    $ws.Cells.Item($r, 9).Value  = $categoriaId
    $ws.Cells.Item($r, 10).Value = $categoria
    $ws.Cells.Item($r, 11).Value = $variedad
    $ws.Cells.Item($r, 12).Value = $nr.Calidad
    $ws.Cells.Item($r, 13).Value = $nr.Volumen
    $ws.Cells.Item($r, 14).Value = $nr.Precio
    $ws.Cells.Item($r, 15).Value = $nr.Precio
    $ws.Cells.Item($r, 16).Value = $nr.Precio
    $ws.Cells.Item($r, 17).Value = $unidad
    $ws.Cells.Item($r, 18).Value = $origen
    $ws.Cells.Item($r, 19).Value = $nr.PrecioKg
    $ws.Cells.Item($r, 20).Value = $kgPorUnidad
}
